# Mise à jour de l'application
# Adds the 5 new "Entrainement" (training) player rows for the J-1 session
# dated 2025-09-12 (serial 45912) at the bottom of the data table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row=477; Player="Malik Boussaid";    Poste="left back";       Temps="01:20:30"; H=6.08; I=0.2;  J=5.87; K=0.18; L=0.02; M=0.01; N=0; O=1; P=3.8;  Q=29.59; R=4.52; S=26; T=1;  U=17; V=13 },
    @{ Row=478; Player="Mattheo Haon";      Poste="right back";      Temps="01:22:04"; H=6.46; I=0.57; J=5.88; K=0.46; L=0.13; M=0;    N=0; O=0; P=4.66; Q=23.57; R=4.73; S=27; T=8;  U=21; V=4  },
    @{ Row=479; Player="Ilan Ihaddadene";   Poste="center midfield"; Temps="01:22:27"; H=6.29; I=0.19; J=6.09; K=0.17; L=0.02; M=0;    N=0; O=1; P=4.49; Q=25.32; R=5.02; S=19; T=3;  U=15; V=0  },
    @{ Row=480; Player="Omar Benyounes";    Poste="center midfield"; Temps="01:22:20"; H=6.38; I=0.38; J=5.99; K=0.31; L=0.08; M=0;    N=0; O=0; P=4.55; Q=23.81; R=4.92; S=32; T=8;  U=16; V=3  },
    @{ Row=481; Player="Ilyes Boughanmi";   Poste="center forward";  Temps="01:21:40"; H=5.34; I=0.39; J=4.94; K=0.33; L=0.07; M=0;    N=0; O=0; P=3.84; Q=22.63; R=5.03; S=18; T=11; U=32; V=5  }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = "Entrainement"

    # Column B (Date) must keep the same date-number format already used by
    # the existing rows (style index reused, no new numFmt created).
    $ws.Range("B476").Copy()
    $ws.Range("B$row").PasteSpecial(-4122)
    $ws.Cells.Item($row, 2).Value = 45912

    $ws.Cells.Item($row, 3).Value = "Global"
    $ws.Cells.Item($row, 4).Value = "J-1"
    $ws.Cells.Item($row, 5).Value = $r.Player
    $ws.Cells.Item($row, 6).Value = $r.Poste
    $ws.Cells.Item($row, 7).Value = $r.Temps
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
    $ws.Cells.Item($row, 21).Value = $r.U
    $ws.Cells.Item($row, 22).Value = $r.V
}

$excel.ActiveWindow.ScrollRow = 459
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F489").Select()
